$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OpsTracker")

# Add new row 36 with the new task entry
$ws.Cells.Item(36, 1).Value = 36
$ws.Cells.Item(36, 2).Value = "Need to take interview of Jayjit Sen on 16th Nov at 2 PM"
$ws.Cells.Item(36, 3).Value = "Debasish"
$ws.Cells.Item(36, 4).Value = "Todo"

# Update the selection on the OpsTracker sheet to match the new active cell
$ws.Activate()
$ws.Range("D37").Select()
